# Applies the "Remove renewables from BAU guaranteed dispatch and set coal
# to bid at its expected capacity factor" edit to the BGDPbES sheet.
#
# The user selected, for each electricity source row, the block of early
# -year cells that still held the "=$B<row>" fill-formula and overtyped it
# with a literal 0 (nuclear, hydro, onshore wind, solar PV, solar thermal,
# biomass, geothermal, offshore wind and municipal solid waste all go to a
# guaranteed-dispatch percentage of 0 -- i.e. renewables, nuclear, and MSW
# are removed from guaranteed dispatch). For hard coal and natural gas
# nonpeaker the 2015-2018 cells already held 0 and only lose their
# formulas (the value doesn't change), matching "bid at its expected
# capacity factor" for the STEO-driven 2020-2022 cells further right,
# which are left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")
$ws.Activate()

# hard coal (row 2) and natural gas nonpeaker (row 3): value already 0,
# just overtype the early fill-formula cells with the literal value.
$ws.Range("C2:F2").Value = 0
$ws.Range("C3:F3").Value = 0

# nuclear, hydro, onshore wind, solar PV, solar thermal, biomass,
# geothermal (rows 4-10): remove guaranteed dispatch entirely (1 -> 0)
$ws.Range("B4:I4").Value = 0
$ws.Range("B5:I5").Value = 0
$ws.Range("B6:I6").Value = 0
$ws.Range("B7:I7").Value = 0
$ws.Range("B8:I8").Value = 0
$ws.Range("B9:I9").Value = 0
$ws.Range("B10:I10").Value = 0

# petroleum (row 11), natural gas peaker (row 12), lignite (row 13):
# value already 0, just overtype the early fill-formula cells.
$ws.Range("C11:I11").Value = 0
$ws.Range("C12:F12").Value = 0
$ws.Range("C13:F13").Value = 0

# offshore wind (row 14): remove guaranteed dispatch entirely (1 -> 0)
$ws.Range("B14:J14").Value = 0

# municipal solid waste (row 17): remove guaranteed dispatch entirely (1 -> 0)
$ws.Range("B17:F17").Value = 0

# Restore the on-screen view/selection to match where the user ended up
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("F14:J14").Select()
